$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of after_row -> source (before) row, derived from the diff.
# Each destination row's D, J, K, L, M, P values are simply copied from
# another row's original values (a row-wise shuffle of the weekly prices).
$mapping = @{
    2  = 15
    3  = 5
    4  = 9
    5  = 2
    6  = 11
    7  = 17
    8  = 16
    9  = 12
    10 = 7
    11 = 4
    12 = 8
    13 = 3
    14 = 13
    15 = 10
    16 = 14
    17 = 6
}

# Snapshot of the original (pre-edit) values for the columns that move:
# D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado), P (Precio $/Kg)
$original = @{
    2  = @{ D = 44964; J = 1000; K = 2000; L = 2500; M = 2250; P = 750 }
    3  = @{ D = 44883; J = 500;  K = 1800; L = 2000; M = 1900; P = 633 }
    4  = @{ D = 44971; J = 1000; K = 2000; L = 2500; M = 2250; P = 750 }
    5  = @{ D = 44827; J = 1200; K = 2000; L = 2500; M = 2250; P = 750 }
    6  = @{ D = 44951; J = 800;  K = 2000; L = 2500; M = 2250; P = 750 }
    7  = @{ D = 44965; J = 1120; K = 2000; L = 2500; M = 2250; P = 750 }
    8  = @{ D = 44953; J = 1000; K = 2000; L = 2500; M = 2250; P = 750 }
    9  = @{ D = 44978; J = 1000; K = 1800; L = 2000; M = 1900; P = 633 }
    10 = @{ D = 44970; J = 800;  K = 2000; L = 2500; M = 2250; P = 750 }
    11 = @{ D = 44992; J = 1040; K = 2000; L = 2500; M = 2250; P = 750 }
    12 = @{ D = 44985; J = 1000; K = 2000; L = 2500; M = 2250; P = 750 }
    13 = @{ D = 44685; J = 400;  K = 1500; L = 2000; M = 1750; P = 583 }
    14 = @{ D = 44881; J = 500;  K = 1900; L = 2000; M = 1950; P = 650 }
    15 = @{ D = 44910; J = 1000; K = 1800; L = 2000; M = 1900; P = 633 }
    16 = @{ D = 44911; J = 700;  K = 1800; L = 2000; M = 1900; P = 633 }
    17 = @{ D = 44848; J = 1000; K = 1500; L = 2000; M = 1750; P = 583 }
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $vals = $original[$srcRow]

    $ws.Cells.Item($destRow, 4).Value  = $vals.D   # D: Fecha
    $ws.Cells.Item($destRow, 10).Value = $vals.J   # J: Volumen
    $ws.Cells.Item($destRow, 11).Value = $vals.K   # K: Precio minimo
    $ws.Cells.Item($destRow, 12).Value = $vals.L   # L: Precio maximo
    $ws.Cells.Item($destRow, 13).Value = $vals.M   # M: Precio promedio ponderado
    $ws.Cells.Item($destRow, 16).Value = $vals.P   # P: Precio $/Kg
}
